$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 47 (Ines Amara): institution/country/year corrected ---
$ws.Range("C47").Value = "Université de Tunis El Manar"
$ws.Range("D47").Value = "Tunisie"
$ws.Range("G47").Value = 1987

# --- Widen column C to fit the longer institution names ---
$ws.Columns.Item(3).ColumnWidth = 48.66

# --- Add new row 56: Abir Smiti ---
$ws.Range("A56").Value = "Abir"
$ws.Range("B56").Value = "Smiti"
$ws.Range("C56").Value = "Université de Tunis"
$ws.Range("D56").Value = "Tunisie"
$ws.Range("E56").Value = "zbDwB7MAAAAJ"
$ws.Range("F56").Value = "F"
$ws.Range("G56").Value = 1985
$ws.Range("H56").Value = "Informatique, Mathématiques et Ingénierie"
$ws.Range("F55").Copy()
$ws.Range("F56").PasteSpecial(-4122)

# --- Add new row 57: Latifa Remadi ---
$ws.Range("A57").Value = "Latifa"
$ws.Range("B57").Value = "Remadi"
$ws.Range("C57").Value = "Foundation for Research and Technology - Hellas"
$ws.Range("D57").Value = "Grèce"
$ws.Range("E57").Value = "c3jenkwAAAAJ"
$ws.Range("F57").Value = "F"
$ws.Range("G57").Value = 1987
$ws.Range("H57").Value = "Médecine, Biologie et Sciences de la Santé"
$ws.Range("F55").Copy()
$ws.Range("F57").PasteSpecial(-4122)

# --- Update the view: active cell below the new last row ---
$ws.Range("H58").Select()

$excel.CutCopyMode = $false
